$wb = $excel.ActiveWorkbook

# --- 1. "sets" sheet: set_id=3 home_points (D4) changes from 15 to 17 ---
$wsSets = $wb.Worksheets.Item("sets")
$wsSets.Range("D4").Value = 17

# --- 2. "rallies" sheet: append two new rally rows (79 and 80) ---
$wsRallies = $wb.Worksheets.Item("rallies")

# Row 79 -> rally_id 78
$wsRallies.Cells.Item(79, 1).Value = 78
$wsRallies.Cells.Item(79, 2).Value = 1
$wsRallies.Cells.Item(79, 3).Value = 3
$wsRallies.Cells.Item(79, 4).Value = 16
$wsRallies.Cells.Item(79, 5).Value = "NOS"
$wsRallies.Cells.Item(79, 6).Value = ""
$wsRallies.Cells.Item(79, 7).Value = 3
$wsRallies.Cells.Item(79, 8).Value = "LINHA"
$wsRallies.Cells.Item(79, 9).Value = "PONTO"
$wsRallies.Cells.Item(79, 10).Value = "NOS"
$wsRallies.Cells.Item(79, 11).Value = 16
$wsRallies.Cells.Item(79, 12).Value = 0
$wsRallies.Cells.Item(79, 13).Value = "1 3 l"
$wsRallies.Cells.Item(79, 14).Value = "FRENTE"
$wsRallies.Cells.Item(79, 15).Value = "FRENTE"
$wsRallies.Cells.Item(79, 16).Value = "FRENTE"

# Row 80 -> rally_id 79
$wsRallies.Cells.Item(80, 1).Value = 79
$wsRallies.Cells.Item(80, 2).Value = 1
$wsRallies.Cells.Item(80, 3).Value = 3
$wsRallies.Cells.Item(80, 4).Value = 17
$wsRallies.Cells.Item(80, 5).Value = "NOS"
$wsRallies.Cells.Item(80, 6).Value = ""
$wsRallies.Cells.Item(80, 7).Value = 5
$wsRallies.Cells.Item(80, 8).Value = "LOB"
$wsRallies.Cells.Item(80, 9).Value = "PONTO"
$wsRallies.Cells.Item(80, 10).Value = "NOS"
$wsRallies.Cells.Item(80, 11).Value = 17
$wsRallies.Cells.Item(80, 12).Value = 0
$wsRallies.Cells.Item(80, 13).Value = "1 5 lob"
$wsRallies.Cells.Item(80, 14).Value = "FRENTE"
$wsRallies.Cells.Item(80, 15).Value = "FRENTE"
$wsRallies.Cells.Item(80, 16).Value = "FRENTE"
